# Auto-generated Excel COM-interop script
# Applies numeric cell updates to match target diff across multiple sheets
# (ALC, ARM, BSM, CRP, CUL, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2823.5
$ws.Range("I40").Value = 1492
$ws.Range("K40").Value = 1492
$ws.Range("M40").Value = -1317
$ws.Range("H74").Value = 4011.2778
$ws.Range("I74").Value = 2707.3572
$ws.Range("J74").Value = 8575
$ws.Range("K74").Value = 2707.3572
$ws.Range("L74").Value = 8575
$ws.Range("M74").Value = -1771.3572
$ws.Range("N74").Value = -10447
$ws.Range("H77").Value = 4011.2778
$ws.Range("I77").Value = 2707.3572
$ws.Range("J77").Value = 8575
$ws.Range("K77").Value = 13536.786
$ws.Range("L77").Value = 42875
$ws.Range("M77").Value = -8856.786
$ws.Range("N77").Value = -52235
$ws.Range("H116").Value = 5695274.5
$ws.Range("J116").Value = 4212.5454
$ws.Range("L116").Value = 4212.5454
$ws.Range("N116").Value = -11096.5454
$ws.Range("H127").Value = 1423.2
$ws.Range("I127").Value = 1176.1428
$ws.Range("J127").Value = 1999.6666
$ws.Range("K127").Value = 3528.4284
$ws.Range("L127").Value = 5998.9998
$ws.Range("M127").Value = 1431.5716
$ws.Range("N127").Value = -15918.9998
$ws.Range("H129").Value = 38463332
$ws.Range("J129").Value = 47620804
$ws.Range("L129").Value = 142862412
$ws.Range("N129").Value = -142872412
$ws.Range("H137").Value = 434931.1
$ws.Range("I137").Value = 793912.3
$ws.Range("J137").Value = 10680.546
$ws.Range("K137").Value = 2381736.9
$ws.Range("L137").Value = 32041.638
$ws.Range("M137").Value = -2379186.9
$ws.Range("N137").Value = -37141.638
$ws.Range("H138").Value = 139153.72
$ws.Range("I138").Value = 879275
$ws.Range("J138").Value = 5396.8677
$ws.Range("K138").Value = 2637825
$ws.Range("L138").Value = 16190.6031
$ws.Range("M138").Value = -2632685
$ws.Range("N138").Value = -26470.6031

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 21898.8
$ws.Range("I6").Value = 14873
$ws.Range("K6").Value = 14873
$ws.Range("M6").Value = -14700
$ws.Range("H43").Value = 21080.75
$ws.Range("J43").Value = 20992.5
$ws.Range("L43").Value = 20992.5
$ws.Range("N43").Value = -21618.5
$ws.Range("H122").Value = 415294.03
$ws.Range("I122").Value = 2704.9
$ws.Range("K122").Value = 8114.700000000001
$ws.Range("M122").Value = -5664.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 50273.457
$ws.Range("I105").Value = 77551.07000000001
$ws.Range("J105").Value = 4810.778
$ws.Range("K105").Value = 77551.07000000001
$ws.Range("L105").Value = 4810.778
$ws.Range("M105").Value = -75804.07000000001
$ws.Range("N105").Value = -8304.778
$ws.Range("H107").Value = 4783.933
$ws.Range("I107").Value = 5605.1665
$ws.Range("K107").Value = 5605.1665
$ws.Range("M107").Value = -3685.1665
$ws.Range("H134").Value = 5370.488
$ws.Range("J134").Value = 2966.3333
$ws.Range("L134").Value = 8898.999899999999
$ws.Range("N134").Value = -13968.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4419
$ws.Range("I16").Value = 5031.6665
$ws.Range("K16").Value = 5031.6665
$ws.Range("M16").Value = -4744.6665
$ws.Range("H31").Value = 2198.0579
$ws.Range("I31").Value = 1307.5333
$ws.Range("J31").Value = 2445.426
$ws.Range("K31").Value = 1307.5333
$ws.Range("L31").Value = 2445.426
$ws.Range("M31").Value = -1012.5333
$ws.Range("N31").Value = -3035.426
$ws.Range("H34").Value = 2198.0579
$ws.Range("I34").Value = 1307.5333
$ws.Range("J34").Value = 2445.426
$ws.Range("K34").Value = 1307.5333
$ws.Range("L34").Value = 2445.426
$ws.Range("M34").Value = -1105.5333
$ws.Range("N34").Value = -2849.426
$ws.Range("H99").Value = 20836332
$ws.Range("J99").Value = 4999.5
$ws.Range("L99").Value = 4999.5
$ws.Range("N99").Value = -7995.5
$ws.Range("H113").Value = 4419
$ws.Range("I113").Value = 5031.6665
$ws.Range("K113").Value = 5031.6665
$ws.Range("M113").Value = -2861.6665
$ws.Range("H126").Value = 20836332
$ws.Range("J126").Value = 4999.5
$ws.Range("L126").Value = 14998.5
$ws.Range("N126").Value = -19938.5
$ws.Range("H132").Value = 21054
$ws.Range("I132").Value = 7986.8965
$ws.Range("K132").Value = 23960.6895
$ws.Range("M132").Value = -21430.6895

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 206.25
$ws.Range("I33").Value = 164.28572
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 985.71432
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -702.71432
$ws.Range("N33").Value = -3566
$ws.Range("H56").Value = 5361
$ws.Range("I56").Value = 5361
$ws.Range("K56").Value = 5361
$ws.Range("M56").Value = -4831
$ws.Range("H68").Value = 13892114
$ws.Range("J68").Value = 19234484
$ws.Range("L68").Value = 57703452
$ws.Range("N68").Value = -57705074
$ws.Range("H71").Value = 13892114
$ws.Range("J71").Value = 19234484
$ws.Range("L71").Value = 173110356
$ws.Range("N71").Value = -173118468
$ws.Range("H97").Value = 41779.668
$ws.Range("I97").Value = 68366.664
$ws.Range("J97").Value = 1899.1666
$ws.Range("K97").Value = 205099.992
$ws.Range("L97").Value = 5697.4998
$ws.Range("M97").Value = -204603.992
$ws.Range("N97").Value = -6689.4998
$ws.Range("H98").Value = 1002.25
$ws.Range("I98").Value = 1003
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 3009
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -1511
$ws.Range("N98").Value = -5996
$ws.Range("H121").Value = 1922332.5
$ws.Range("I121").Value = 2110178.8
$ws.Range("K121").Value = 6330536.399999999
$ws.Range("M121").Value = -6329226.399999999
$ws.Range("H129").Value = 37039868
$ws.Range("J129").Value = 41669664
$ws.Range("L129").Value = 125008992
$ws.Range("N129").Value = -125018992
$ws.Range("H137").Value = 11695.3
$ws.Range("J137").Value = 14994.167
$ws.Range("L137").Value = 44982.501
$ws.Range("N137").Value = -55182.501
$ws.Range("H140").Value = 3435.7778
$ws.Range("I140").Value = 3435.7778
$ws.Range("K140").Value = 10307.3334
$ws.Range("M140").Value = -5127.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1933
$ws.Range("I16").Value = 1933
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1933
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1763
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 3276.182
$ws.Range("I46").Value = 683.3333
$ws.Range("K46").Value = 683.3333
$ws.Range("M46").Value = -495.3333
$ws.Range("H61").Value = 3218.7368
$ws.Range("I61").Value = 3147
$ws.Range("K61").Value = 3147
$ws.Range("M61").Value = -2945
$ws.Range("H68").Value = 6044.222
$ws.Range("I68").Value = 3300
$ws.Range("J68").Value = 6828.2856
$ws.Range("K68").Value = 3300
$ws.Range("L68").Value = 6828.2856
$ws.Range("M68").Value = -2551
$ws.Range("N68").Value = -8326.285599999999
$ws.Range("H71").Value = 6044.222
$ws.Range("I71").Value = 3300
$ws.Range("J71").Value = 6828.2856
$ws.Range("K71").Value = 16500
$ws.Range("L71").Value = 34141.428
$ws.Range("M71").Value = -12756
$ws.Range("N71").Value = -41629.428
$ws.Range("H100").Value = 15832.333
$ws.Range("I100").Value = 28000
$ws.Range("J100").Value = 3664.6667
$ws.Range("K100").Value = 28000
$ws.Range("L100").Value = 3664.6667
$ws.Range("M100").Value = -27459
$ws.Range("N100").Value = -4746.6667
$ws.Range("H113").Value = 3218.7368
$ws.Range("I113").Value = 3147
$ws.Range("K113").Value = 3147
$ws.Range("M113").Value = -977
$ws.Range("H122").Value = 4924.615
$ws.Range("I122").Value = 4844
$ws.Range("K122").Value = 14532
$ws.Range("M122").Value = -12082
$ws.Range("H131").Value = 67500
$ws.Range("J131").Value = 67500
$ws.Range("L131").Value = 67500
$ws.Range("N131").Value = -77580
$ws.Range("H132").Value = 673896.6
$ws.Range("I132").Value = 986628.6
$ws.Range("J132").Value = 3756.7144
$ws.Range("K132").Value = 2959885.8
$ws.Range("L132").Value = 11270.1432
$ws.Range("M132").Value = -2957355.8
$ws.Range("N132").Value = -16330.1432
$ws.Range("H136").Value = 6573.9414
$ws.Range("I136").Value = 5483
$ws.Range("J136").Value = 11665
$ws.Range("K136").Value = 16449
$ws.Range("L136").Value = 34995
$ws.Range("M136").Value = -13899
$ws.Range("N136").Value = -40095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3792.8462
$ws.Range("I113").Value = 2400.5789
$ws.Range("K113").Value = 7201.736699999999
$ws.Range("M113").Value = -5031.736699999999
